# fix validate excel + add payment success
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# A2: flight number QH3456 -> QH4444
$ws.Range("A2").Value = "QH4444"

# B2: stays QH1111 (unchanged, but ensure it's set explicitly)
$ws.Range("B2").Value = "QH1111"

# C2: departure time moves from 08:00 to 20:00 on the same day (45602.333333333336 -> 45602.833333333336)
$ws.Range("C2").Value = 45602.833333333336

# Update the active selection to E8
$ws.Range("E8").Select()
